$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 406.1
$ws.Range("I6").Value = 65.25
$ws.Range("J6").Value = 633.3333
$ws.Range("K6").Value = 195.75
$ws.Range("L6").Value = 1899.9999
$ws.Range("M6").Value = -83.75
$ws.Range("N6").Value = -2123.9999

$ws.Range("H8").Value = 500.2
$ws.Range("I8").Value = 500.2
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1500.6
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -1361.6
$ws.Range("N8").ClearContents()

$ws.Range("H64").Value = 45457470
$ws.Range("I64").Value = 71430824
$ws.Range("J64").Value = 4107.5
$ws.Range("K64").Value = 71430824
$ws.Range("L64").Value = 4107.5
$ws.Range("M64").Value = -71430576
$ws.Range("N64").Value = -4603.5

$ws.Range("H67").Value = 45457470
$ws.Range("I67").Value = 71430824
$ws.Range("J67").Value = 4107.5
$ws.Range("K67").Value = 71430824
$ws.Range("L67").Value = 4107.5
$ws.Range("M67").Value = -71429966
$ws.Range("N67").Value = -5823.5

$ws.Range("H74").Value = 2710
$ws.Range("I74").Value = 2041.1852
$ws.Range("J74").Value = 3495.1304
$ws.Range("K74").Value = 2041.1852
$ws.Range("L74").Value = 3495.1304
$ws.Range("M74").Value = -1105.1852
$ws.Range("N74").Value = -5367.1304

$ws.Range("H77").Value = 2710
$ws.Range("I77").Value = 2041.1852
$ws.Range("J77").Value = 3495.1304
$ws.Range("K77").Value = 10205.926
$ws.Range("L77").Value = 17475.652
$ws.Range("M77").Value = -5525.925999999999
$ws.Range("N77").Value = -26835.652

$ws.Range("H87").Value = 36925
$ws.Range("J87").Value = 36925
$ws.Range("L87").Value = 36925
$ws.Range("N87").Value = -39421

$ws.Range("H90").Value = 36925
$ws.Range("J90").Value = 36925
$ws.Range("L90").Value = 110775
$ws.Range("N90").Value = -123255

$ws.Range("H98").Value = 2221.5
$ws.Range("I98").Value = 1782.2307
$ws.Range("K98").Value = 1782.2307
$ws.Range("M98").Value = -284.2307000000001

$ws.Range("H122").Value = 2221.5
$ws.Range("I122").Value = 1782.2307
$ws.Range("K122").Value = 5346.6921
$ws.Range("M122").Value = -2896.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H61").Value = 229049.95
$ws.Range("I61").Value = 2064.05
$ws.Range("J61").Value = 418204.88
$ws.Range("K61").Value = 2064.05
$ws.Range("L61").Value = 418204.88
$ws.Range("M61").Value = -1852.05
$ws.Range("N61").Value = -418628.88

$ws.Range("H97").Value = 691.3421
$ws.Range("I97").Value = 596.7406999999999
$ws.Range("J97").Value = 923.5454999999999
$ws.Range("K97").Value = 596.7406999999999
$ws.Range("L97").Value = 923.5454999999999
$ws.Range("M97").Value = -100.7406999999999
$ws.Range("N97").Value = -1915.5455

$ws.Range("H102").Value = 1204.2858
$ws.Range("I102").Value = 1066
$ws.Range("J102").Value = 1550
$ws.Range("K102").Value = 1066
$ws.Range("L102").Value = 1550
$ws.Range("M102").Value = 556
$ws.Range("N102").Value = -4794

$ws.Range("H122").Value = 1878.6923
$ws.Range("I122").Value = 1985.3334
$ws.Range("J122").Value = 1787.2858
$ws.Range("K122").Value = 5956.0002
$ws.Range("L122").Value = 5361.857400000001
$ws.Range("M122").Value = -3506.0002
$ws.Range("N122").Value = -10261.8574

$ws.Range("H136").Value = 229049.95
$ws.Range("I136").Value = 2064.05
$ws.Range("J136").Value = 418204.88
$ws.Range("K136").Value = 6192.150000000001
$ws.Range("L136").Value = 1254614.64
$ws.Range("M136").Value = -3642.150000000001
$ws.Range("N136").Value = -1259714.64

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 874.19354
$ws.Range("I94").Value = 609.375
$ws.Range("J94").Value = 1156.6666
$ws.Range("K94").Value = 609.375
$ws.Range("L94").Value = 1156.6666
$ws.Range("M94").Value = -158.375
$ws.Range("N94").Value = -2058.6666

$ws.Range("H105").Value = 1588.3182
$ws.Range("I105").Value = 1144.2858
$ws.Range("J105").Value = 2365.375
$ws.Range("K105").Value = 1144.2858
$ws.Range("L105").Value = 2365.375
$ws.Range("M105").Value = 602.7141999999999
$ws.Range("N105").Value = -5859.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 7500
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 7500
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 7500
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -8452

$ws.Range("H99").Value = 1799.091
$ws.Range("I99").Value = 1664
$ws.Range("J99").Value = 2407
$ws.Range("K99").Value = 1664
$ws.Range("L99").Value = 2407
$ws.Range("M99").Value = -166
$ws.Range("N99").Value = -5403

$ws.Range("H126").Value = 1799.091
$ws.Range("I126").Value = 1664
$ws.Range("J126").Value = 2407
$ws.Range("K126").Value = 4992
$ws.Range("L126").Value = 7221
$ws.Range("M126").Value = -2522
$ws.Range("N126").Value = -12161

$ws.Range("H134").Value = 4333.6665
$ws.Range("I134").Value = 4364
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 13092
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -10557
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4197.1724
$ws.Range("I3").Value = 2794.577
$ws.Range("K3").Value = 8383.731
$ws.Range("M3").Value = -8271.731

$ws.Range("H59").Value = 1200
$ws.Range("J59").Value = 1200
$ws.Range("L59").Value = 3600
$ws.Range("N59").Value = -4680

$ws.Range("H131").Value = 918.85187
$ws.Range("J131").Value = 963.76
$ws.Range("L131").Value = 2891.28
$ws.Range("N131").Value = -12971.28

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2431.2307
$ws.Range("I122").Value = 949.625
$ws.Range("K122").Value = 2848.875
$ws.Range("M122").Value = -398.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2654.7222
$ws.Range("I40").Value = 2705.9092
$ws.Range("J40").Value = 2574.2856
$ws.Range("K40").Value = 2705.9092
$ws.Range("L40").Value = 2574.2856
$ws.Range("M40").Value = -2569.9092
$ws.Range("N40").Value = -2846.2856

$ws.Range("H93").Value = 2173
$ws.Range("I93").Value = 1549.75
$ws.Range("J93").Value = 3004
$ws.Range("K93").Value = 1549.75
$ws.Range("L93").Value = 3004
$ws.Range("M93").Value = -301.75
$ws.Range("N93").Value = -5500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 872.65
$ws.Range("I126").Value = 723.3125
$ws.Range("J126").Value = 1470
$ws.Range("K126").Value = 2169.9375
$ws.Range("L126").Value = 4410
$ws.Range("M126").Value = 300.0625
$ws.Range("N126").Value = -9350
